# PowerShell Excel COM-interop script
# Applies updated TPM-derived values to Sheet1 of the workbook (Col1a2-Itgb3 LR-pairs)
# Each data line is: row,col,newValue
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellData = @"
2,7,12.73908466666667
2,8,38.217254
2,9,0.005953388968763418
2,10,0.006105597140986208
2,13,6.712486666666666
2,14,20.13746
2,15,0.6330487633990675
2,16,0.6414503882251803
2,17,85.51093597053779
2,18,769.59842373484
2,19,0.003768785524709331
2,20,0.003916437656432154
3,7,12.73908466666667
3,8,38.217254
3,9,0.005953388968763418
3,10,0.006105597140986208
3,15,0.290741083484562
3,16,0.2945997080427384
3,17,39.27271264280268
3,18,353.4544137852241
3,19,0.001730894759183315
3,20,0.001798707135161115
4,7,12.73908466666667
4,8,38.217254
4,9,0.005953388968763418
4,10,0.006105597140986208
4,13,0.2495096666666667
4,14,0.748529
4,15,0.02353103905946135
4,16,0.02384333563656022
4,17,3.178524768818445
4,18,28.606722919366
4,19,0.0001400894283601383
4,20,0.0001455778018941567
5,7,12.73908466666667
5,8,38.217254
5,9,0.005953388968763418
5,10,0.006105597140986208
5,13,0.4166465
5,14,0.8332930000000001
5,15,0.03929356804674715
5,16,0.02654337331298611
5,17,5.307695039570334
5,18,31.846170237422
5,19,0.0002339298945528592
5,20,0.0001620631442118976
6,7,12.73908466666667
6,8,38.217254
6,9,0.005953388968763418
6,10,0.006105597140986208
6,11,2
6,12,0.6666666666666666
6,13,0.1419326666666667
6,14,0.425798
6,15,0.01338554601016197
6,16,0.01356319478253491
6,17,1.808092257632445
6,18,16.272830318692
6,19,7.968936195777343E-05
6,20,8.281140328688417E-05
7,9,0.9182810852447438
7,10,0.9417584502053091
7,13,6.712486666666666
7,14,20.13746
7,15,0.6330487633990675
7,16,0.6414503882251803
7,17,13189.64299079375
7,18,118706.7869171437
7,19,0.5813167054669387
7,20,0.6040913234985397
8,9,0.9182810852447438
8,10,0.9417584502053091
8,15,0.290741083484562
8,16,0.2945997080427384
8,19,0.2669820376674362
8,20,0.2774417644772658
9,9,0.9182810852447438
9,10,0.9417584502053091
9,13,0.2495096666666667
9,14,0.748529
9,15,0.02353103905946135
9,16,0.02384333563656022
9,17,490.2718753137612
9,18,4412.446877823851
9,19,0.02160810808445863
9,20,0.02245466281681198
10,9,0.9182810852447438
10,10,0.9417584502053091
10,13,0.4166465
10,14,0.8332930000000001
10,15,0.03929356804674715
10,16,0.02654337331298611
10,17,818.6859596538612
10,18,4912.115757923168
10,19,0.03608254030910516
10,20,0.02499744611445876
11,9,0.9182810852447438
11,10,0.9417584502053091
11,11,2
11,12,0.6666666666666666
11,13,0.1419326666666667
11,14,0.425798
11,15,0.01338554601016197
11,16,0.01356319478253491
11,17,278.8893736446403
11,18,2510.004362801762
11,19,0.01229169371680498
11,20,0.01277325329823281
12,7,1.091866333333334
12,8,3.275599000000001
12,9,0.0005102646818291153
12,10,0.0005233104369407934
12,13,6.712486666666666
12,14,20.13746
12,15,0.6330487633990675
12,16,0.6414503882251803
12,17,7.329138204282224
12,18,65.96224383854
12,19,0.00032302242583814
12,20,0.0003356776829379606
13,7,1.091866333333334
13,8,3.275599000000001
13,9,0.0005102646818291153
13,10,0.0005233104369407934
13,15,0.290741083484562
13,16,0.2945997080427384
13,17,3.366062309449334
13,18,30.29456078504401
13,19,0.0001483549064589023
13,20,0.0001541671019384756
14,7,1.091866333333334
14,8,3.275599000000001
14,9,0.0005102646818291153
14,10,0.0005233104369407934
14,13,0.2495096666666667
14,14,0.748529
14,15,0.02353103905946135
14,16,0.02384333563656022
14,17,0.2724312048745556
14,18,2.451880843871
14,19,1.200705815878453E-05
14,20,1.247746639009432E-05
15,7,1.091866333333334
15,8,3.275599000000001
15,9,0.0005102646818291153
15,10,0.0005233104369407934
15,13,0.4166465
15,14,0.8332930000000001
15,15,0.03929356804674715
15,16,0.02654337331298611
15,17,0.4549222862511668
15,18,2.729533717507001
15,19,2.005011999730413E-05
15,20,1.389042428630136E-05
16,7,1.091866333333334
16,8,3.275599000000001
16,9,0.0005102646818291153
16,10,0.0005233104369407934
16,11,2
16,12,0.6666666666666666
16,13,0.1419326666666667
16,14,0.425798
16,15,0.01338554601016197
16,16,0.01356319478253491
16,17,0.1549715003335556
16,18,1.394743503002
16,19,6.830171375984279E-06
16,20,7.097761387961431E-06
17,7,160.0313415
17,8,320.062683
17,9,0.0747878554913321
17,10,0.05113328661083746
17,13,6.712486666666666
17,14,20.13746
17,15,0.6330487633990675
17,16,0.6414503882251803
17,17,1074.20824606753
17,18,6445.249476405179
17,19,0.04734435943605594
17,20,0.03279946654775111
18,7,160.0313415
18,8,320.062683
18,9,0.0747878554913321
18,10,0.05113328661083746
18,15,0.290741083484562
18,16,0.2945997080427384
18,17,493.3529412059581
18,18,2960.117647235748
18,19,0.02174390213703674
18,20,0.01506385130681838
19,7,160.0313415
19,8,320.062683
19,9,0.0747878554913321
19,10,0.05113328661083746
19,13,0.2495096666666667
19,14,0.748529
19,15,0.02353103905946135
19,16,0.02384333563656022
19,17,39.92936667388449
19,18,239.576200043307
19,19,0.001759835948739887
19,20,0.001219188114862629
20,7,160.0313415
20,8,320.062683
20,9,0.0747878554913321
20,10,0.05113328661083746
20,13,0.4166465
20,14,0.8332930000000001
20,15,0.03929356804674715
20,16,0.02654337331298611
20,17,66.67649832627976
20,18,266.705993305119
20,19,0.002938681688818951
20,20,0.001357249915231373
21,7,160.0313415
21,8,320.062683
21,9,0.0747878554913321
21,10,0.05113328661083746
21,11,2
21,12,0.6666666666666666
21,13,0.1419326666666667
21,14,0.425798
21,15,0.01338554601016197
21,16,0.01356319478253491
21,17,22.713675049339
21,18,136.282050296034
21,19,0.00100107628068057
21,20,0.0006935307261739727
22,7,1.000156333333333
22,8,3.000469
22,9,0.0004674056133315229
22,10,0.0004793556059265206
22,13,6.712486666666666
22,14,20.13746
22,15,0.6330487633990675
22,16,0.6414503882251803
22,17,6.713536052082221
22,18,60.42182446873999
22,19,0.0002958905455253033
22,20,0.0003074828395194832
23,7,1.000156333333333
23,8,3.000469
23,9,0.0004674056133315229
23,10,0.0004793556059265206
23,15,0.290741083484562
23,16,0.2945997080427384
23,17,3.083333952529333
23,18,27.750005572764
23,19,0.0001358940144467732
23,20,0.0001412180215546029
24,7,1.000156333333333
24,8,3.000469
24,9,0.0004674056133315229
24,10,0.0004793556059265206
24,13,0.2495096666666667
24,14,0.748529
24,15,0.02353103905946135
24,16,0.02384333563656022
24,17,0.2495486733445555
24,18,2.245938060101
24,19,1.099853974391556E-05
24,20,1.142943660137273E-05
25,7,1.000156333333333
25,8,3.000469
25,9,0.0004674056133315229
25,10,0.0004793556059265206
25,13,0.4166465
25,14,0.8332930000000001
25,15,0.03929356804674715
25,16,0.02654337331298611
25,17,0.4167116357361667
25,18,2.500269814417
25,19,1.836603427287378E-05
25,20,1.272371479778029E-05
26,7,1.000156333333333
26,8,3.000469
26,9,0.0004674056133315229
26,10,0.0004793556059265206
26,11,2
26,12,0.6666666666666666
26,13,0.1419326666666667
26,14,0.425798
26,15,0.01338554601016197
26,16,0.01356319478253491
26,17,0.1419548554735555
26,18,1.277593699262
26,19,6.256479342657074E-06
26,20,6.501593453281442E-06
"@

$lines = $cellData -split "`n" | Where-Object { $_.Trim() -ne "" }
foreach ($line in $lines) {
    $parts = $line.Trim() -split ","
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $v = [double]$parts[2]
    $ws.Cells.Item($r, $c).Value = $v
}
